$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 10 (B/C): Objetivos: value becomes the responsible teacher string
$ws.Range("B10").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C10").Value = "5840730 - Antonio Jefferson da Silva Machado"

# Row 13 (B/C): Programa resumido: value becomes "01/01/2023".
# Copy/PasteSpecial (values) from a cell that already holds this exact text
# (B8/C8, under "Ativação:") so Excel keeps it typed as text instead of
# re-parsing the date-looking string into a numeric date serial.
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C13").PasteSpecial(-4163)

# Row 15 (B/C): Programa: value becomes the responsible teacher string
$ws.Range("B15").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C15").Value = "5840730 - Antonio Jefferson da Silva Machado"

# Row 18 (B/C): Método: value becomes the second teacher string
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
